$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44504
$ws.Cells.Item(2, 10).Value = 200
$ws.Cells.Item(3, 4).Value = 44305
$ws.Cells.Item(3, 10).Value = 100
$ws.Cells.Item(3, 11).Value = 2500
$ws.Cells.Item(3, 12).Value = 2500
$ws.Cells.Item(3, 13).Value = 2500
$ws.Cells.Item(3, 14).Value = '$/unidad'
$ws.Cells.Item(3, 16).Value = 2500
$ws.Cells.Item(4, 4).Value = 44477
$ws.Cells.Item(4, 8).Value = 'Sin especificar'
$ws.Cells.Item(4, 9).Value = 'Primera'
$ws.Cells.Item(4, 10).Value = 80
$ws.Cells.Item(4, 11).Value = 800
$ws.Cells.Item(4, 12).Value = 800
$ws.Cells.Item(4, 13).Value = 800
$ws.Cells.Item(4, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(4, 15).Value = 'Perú'
$ws.Cells.Item(4, 16).Value = 800
$ws.Cells.Item(5, 4).Value = 44312
$ws.Cells.Item(5, 8).Value = 'Sin especificar'
$ws.Cells.Item(5, 10).Value = 180
$ws.Cells.Item(5, 11).Value = 2500
$ws.Cells.Item(5, 12).Value = 2500
$ws.Cells.Item(5, 13).Value = 2500
$ws.Cells.Item(5, 15).Value = 'Perú'
$ws.Cells.Item(5, 16).Value = 2500
$ws.Cells.Item(6, 4).Value = 44497
$ws.Cells.Item(6, 8).Value = 'Sin especificar'
$ws.Cells.Item(6, 9).Value = 'Primera'
$ws.Cells.Item(6, 10).Value = 250
$ws.Cells.Item(6, 11).Value = 800
$ws.Cells.Item(6, 12).Value = 800
$ws.Cells.Item(6, 13).Value = 800
$ws.Cells.Item(6, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(6, 15).Value = 'Perú'
$ws.Cells.Item(6, 16).Value = 800
$ws.Cells.Item(7, 4).Value = 44495
$ws.Cells.Item(7, 8).Value = 'Sin especificar'
$ws.Cells.Item(7, 9).Value = 'Primera'
$ws.Cells.Item(7, 10).Value = 200
$ws.Cells.Item(7, 11).Value = 800
$ws.Cells.Item(7, 12).Value = 800
$ws.Cells.Item(7, 13).Value = 800
$ws.Cells.Item(7, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(7, 15).Value = 'Perú'
$ws.Cells.Item(7, 16).Value = 800
$ws.Cells.Item(10, 4).Value = 44167
$ws.Cells.Item(10, 9).Value = 'Primera'
$ws.Cells.Item(10, 10).Value = 400
$ws.Cells.Item(10, 11).Value = 5000
$ws.Cells.Item(10, 12).Value = 5000
$ws.Cells.Item(10, 13).Value = 5000
$ws.Cells.Item(10, 16).Value = 5000
$ws.Cells.Item(11, 4).Value = 44167
$ws.Cells.Item(11, 9).Value = 'Segunda'
$ws.Cells.Item(11, 10).Value = 560
$ws.Cells.Item(12, 4).Value = 44167
$ws.Cells.Item(12, 9).Value = 'Tercera'
$ws.Cells.Item(12, 10).Value = 450
$ws.Cells.Item(12, 11).Value = 2000
$ws.Cells.Item(12, 12).Value = 2000
$ws.Cells.Item(12, 13).Value = 2000
$ws.Cells.Item(12, 14).Value = '$/unidad'
$ws.Cells.Item(12, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(12, 16).Value = 2000
$ws.Cells.Item(13, 4).Value = 44488
$ws.Cells.Item(13, 10).Value = 150
$ws.Cells.Item(13, 11).Value = 800
$ws.Cells.Item(13, 12).Value = 800
$ws.Cells.Item(13, 13).Value = 800
$ws.Cells.Item(13, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(13, 16).Value = 800
$ws.Cells.Item(14, 4).Value = 44491
$ws.Cells.Item(14, 10).Value = 150
$ws.Cells.Item(14, 11).Value = 800
$ws.Cells.Item(14, 12).Value = 800
$ws.Cells.Item(14, 13).Value = 800
$ws.Cells.Item(14, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(14, 15).Value = 'Perú'
$ws.Cells.Item(14, 16).Value = 800
$ws.Cells.Item(15, 4).Value = 44483
$ws.Cells.Item(15, 9).Value = 'Primera'
$ws.Cells.Item(15, 10).Value = 120
$ws.Cells.Item(15, 11).Value = 800
$ws.Cells.Item(15, 12).Value = 800
$ws.Cells.Item(15, 13).Value = 800
$ws.Cells.Item(15, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(15, 15).Value = 'Perú'
$ws.Cells.Item(15, 16).Value = 800
$ws.Cells.Item(16, 4).Value = 44194
$ws.Cells.Item(16, 9).Value = 'Extra'
$ws.Cells.Item(16, 10).Value = 120
$ws.Cells.Item(16, 11).Value = 3500
$ws.Cells.Item(16, 12).Value = 3500
$ws.Cells.Item(16, 13).Value = 3500
$ws.Cells.Item(16, 16).Value = 3500
$ws.Cells.Item(17, 4).Value = 44194
$ws.Cells.Item(17, 11).Value = 3000
$ws.Cells.Item(17, 12).Value = 3000
$ws.Cells.Item(17, 13).Value = 3000
$ws.Cells.Item(17, 14).Value = '$/unidad'
$ws.Cells.Item(17, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(17, 16).Value = 3000
$ws.Cells.Item(18, 4).Value = 44223
$ws.Cells.Item(18, 8).Value = 'Americana O Klondike'
$ws.Cells.Item(18, 9).Value = 'Extra'
$ws.Cells.Item(18, 10).Value = 340
$ws.Cells.Item(18, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(19, 4).Value = 44223
$ws.Cells.Item(19, 8).Value = 'Americana O Klondike'
$ws.Cells.Item(19, 10).Value = 400
$ws.Cells.Item(19, 11).Value = 2000
$ws.Cells.Item(19, 12).Value = 2000
$ws.Cells.Item(19, 13).Value = 2000
$ws.Cells.Item(19, 14).Value = '$/unidad'
$ws.Cells.Item(19, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(19, 16).Value = 2000
$ws.Cells.Item(20, 4).Value = 44223
$ws.Cells.Item(20, 8).Value = 'Americana O Klondike'
$ws.Cells.Item(20, 9).Value = 'Segunda'
$ws.Cells.Item(20, 10).Value = 300
$ws.Cells.Item(20, 11).Value = 1500
$ws.Cells.Item(20, 12).Value = 1500
$ws.Cells.Item(20, 13).Value = 1500
$ws.Cells.Item(20, 14).Value = '$/unidad'
$ws.Cells.Item(20, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(20, 16).Value = 1500
$ws.Cells.Item(21, 4).Value = 44223
$ws.Cells.Item(21, 8).Value = 'Americana O Klondike'
$ws.Cells.Item(21, 9).Value = 'Tercera'
$ws.Cells.Item(21, 10).Value = 160
$ws.Cells.Item(21, 11).Value = 1000
$ws.Cells.Item(21, 12).Value = 1000
$ws.Cells.Item(21, 13).Value = 1000
$ws.Cells.Item(21, 14).Value = '$/unidad'
$ws.Cells.Item(21, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(21, 16).Value = 1000
$ws.Cells.Item(22, 4).Value = 44510
$ws.Cells.Item(22, 10).Value = 250
